$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Range("B20").Value = 6920351
$ws.Range("E20").Value = "Hamburg SV II"
$ws.Range("F20").Value = "SV DrochtersenAssel"
$ws.Range("K20").Value = "D"
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 3
$ws.Range("M20").Value = 3.6
$ws.Range("N20").Value = 2
$ws.Range("O20").Value = 5
$ws.Range("P20").Value = 4
$ws.Range("Q20").Value = 1.615
$ws.Range("R20").Value = 1
$ws.Range("S20").Value = 1.775
$ws.Range("T20").Value = 2.025
$ws.Range("U20").Value = 2.75
$ws.Range("V20").Value = 1.775
$ws.Range("W20").Value = 2.025
$ws.Range("X20").Value = -1
$ws.Range("Y20").Value = 3
$ws.Range("Z20").Value = -1
$ws.Range("AA20").Value = 0.7749999999999999
$ws.Range("AB20").Value = -1
$ws.Range("AC20").Value = -1
$ws.Range("AD20").Value = 1.025

# Row 21
$ws.Range("B21").Value = 6920350
$ws.Range("E21").Value = "Phonix Lubeck"
$ws.Range("F21").Value = "SSV Jeddeloh"
$ws.Range("K21").Value = "H"
$ws.Range("G21").Value = 7
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 3
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 2.25
$ws.Range("M21").Value = 3.75
$ws.Range("N21").Value = 2.5
$ws.Range("O21").Value = 2.05
$ws.Range("P21").Value = 4.2
$ws.Range("Q21").Value = 2.75
$ws.Range("R21").Value = -0.25
$ws.Range("S21").Value = 1.875
$ws.Range("T21").Value = 1.925
$ws.Range("U21").Value = 3
$ws.Range("V21").Value = 1.875
$ws.Range("W21").Value = 1.925
$ws.Range("X21").Value = 1.05
$ws.Range("Y21").Value = -1
$ws.Range("Z21").Value = -1
$ws.Range("AA21").Value = 0.875
$ws.Range("AB21").Value = -1
$ws.Range("AC21").Value = 0.875
$ws.Range("AD21").Value = -1

# Row 293
$ws.Range("B293").Value = 7093876
$ws.Range("E293").Value = "Holstein Kiel II"
$ws.Range("F293").Value = "FC Kilia Kiel"
$ws.Range("K293").Value = "D"
$ws.Range("G293").Value = 3
$ws.Range("H293").Value = 3
$ws.Range("I293").Value = 1
$ws.Range("J293").Value = 2
$ws.Range("L293").Value = 1.285
$ws.Range("M293").Value = 5.5
$ws.Range("N293").Value = 6.5
$ws.Range("O293").Value = 1.25
$ws.Range("P293").Value = 6.25
$ws.Range("Q293").Value = 6.5
$ws.Range("R293").Value = -2
$ws.Range("S293").Value = 2.025
$ws.Range("T293").Value = 1.825
$ws.Range("U293").Value = 4.25
$ws.Range("V293").Value = 1.975
$ws.Range("W293").Value = 1.875
$ws.Range("X293").Value = -1
$ws.Range("Y293").Value = 5.25
$ws.Range("Z293").Value = -1
$ws.Range("AA293").Value = -1
$ws.Range("AB293").Value = 0.825
$ws.Range("AC293").Value = 0.9750000000000001
$ws.Range("AD293").Value = -1

# Row 294
$ws.Range("B294").Value = 7096258
$ws.Range("E294").Value = "Phonix Lubeck"
$ws.Range("F294").Value = "Bremer SV"
$ws.Range("K294").Value = "A"
$ws.Range("G294").Value = 3
$ws.Range("H294").Value = 4
$ws.Range("I294").Value = 1
$ws.Range("J294").Value = 2
$ws.Range("L294").Value = 1.333
$ws.Range("M294").Value = 5
$ws.Range("N294").Value = 6
$ws.Range("O294").Value = 1.533
$ws.Range("P294").Value = 4.333
$ws.Range("Q294").Value = 4.2
$ws.Range("R294").Value = -1
$ws.Range("S294").Value = 1.875
$ws.Range("T294").Value = 1.975
$ws.Range("U294").Value = 3.5
$ws.Range("V294").Value = 1.975
$ws.Range("W294").Value = 1.875
$ws.Range("X294").Value = -1
$ws.Range("Y294").Value = -1
$ws.Range("Z294").Value = 3.2
$ws.Range("AA294").Value = -1
$ws.Range("AB294").Value = 0.9750000000000001
$ws.Range("AC294").Value = 0.9750000000000001
$ws.Range("AD294").Value = -1

# Row 295
$ws.Range("B295").Value = 7096259
$ws.Range("E295").Value = "SC Weiche Flensburg 08"
$ws.Range("F295").Value = "Eintracht Norderstedt"
$ws.Range("K295").Value = "A"
$ws.Range("G295").Value = 0
$ws.Range("H295").Value = 2
$ws.Range("I295").Value = 0
$ws.Range("J295").Value = 0
$ws.Range("L295").Value = 1.95
$ws.Range("M295").Value = 4
$ws.Range("N295").Value = 2.875
$ws.Range("O295").Value = 2.15
$ws.Range("P295").Value = 3.9
$ws.Range("Q295").Value = 2.6
$ws.Range("R295").Value = -0.25
$ws.Range("S295").Value = 1.975
$ws.Range("T295").Value = 1.875
$ws.Range("U295").Value = 3
$ws.Range("V295").Value = 1.875
$ws.Range("W295").Value = 1.975
$ws.Range("X295").Value = -1
$ws.Range("Y295").Value = -1
$ws.Range("Z295").Value = 1.6
$ws.Range("AA295").Value = -1
$ws.Range("AB295").Value = 0.875
$ws.Range("AC295").Value = -1
$ws.Range("AD295").Value = 0.9750000000000001

# Row 296
$ws.Range("B296").Value = 7123930
$ws.Range("E296").Value = "TuS BlauWeiss Lohne"
$ws.Range("F296").Value = "VfB Oldenburg"
$ws.Range("K296").Value = "A"
$ws.Range("G296").Value = 0
$ws.Range("H296").Value = 2
$ws.Range("I296").Value = 0
$ws.Range("J296").Value = 1
$ws.Range("L296").Value = 2.25
$ws.Range("M296").Value = 3.75
$ws.Range("N296").Value = 2.5
$ws.Range("O296").Value = 2.55
$ws.Range("P296").Value = 3.8
$ws.Range("Q296").Value = 2.15
$ws.Range("R296").Value = 0.25
$ws.Range("S296").Value = 1.8
$ws.Range("T296").Value = 2
$ws.Range("U296").Value = 3.25
$ws.Range("V296").Value = 1.975
$ws.Range("W296").Value = 1.825
$ws.Range("X296").Value = -1
$ws.Range("Y296").Value = -1
$ws.Range("Z296").Value = 1.15
$ws.Range("AA296").Value = -1
$ws.Range("AB296").Value = 1
$ws.Range("AC296").Value = -1
$ws.Range("AD296").Value = 0.825

# Row 297
$ws.Range("B297").Value = 7093875
$ws.Range("E297").Value = "Eimsbutteler TV"
$ws.Range("F297").Value = "SSV Jeddeloh"
$ws.Range("K297").Value = "A"
$ws.Range("G297").Value = 0
$ws.Range("H297").Value = 2
$ws.Range("I297").Value = 0
$ws.Range("J297").Value = 1
$ws.Range("L297").Value = 2.5
$ws.Range("M297").Value = 4
$ws.Range("N297").Value = 2.15
$ws.Range("O297").Value = 3.7
$ws.Range("P297").Value = 4.75
$ws.Range("Q297").Value = 1.571
$ws.Range("R297").Value = 1
$ws.Range("S297").Value = 1.825
$ws.Range("T297").Value = 1.975
$ws.Range("U297").Value = 3.5
$ws.Range("V297").Value = 1.85
$ws.Range("W297").Value = 1.95
$ws.Range("X297").Value = -1
$ws.Range("Y297").Value = -1
$ws.Range("Z297").Value = 0.571
$ws.Range("AA297").Value = -1
$ws.Range("AB297").Value = 0.9750000000000001
$ws.Range("AC297").Value = -1
$ws.Range("AD297").Value = 0.95
